$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.391.31"
$ws.Range("E2").Value = "  -1.90%  "

$ws.Range("D3").Value = "2.432.68"
$ws.Range("E3").Value = "  -1.65%  "

$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").Value = "'571.78"
$ws.Range("E5").Value = "  -0.85%  "

$ws.Range("D6").Value = "'143.57"
$ws.Range("E6").Value = "  -3.41%  "

$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("E8").Value = "  -2.15%  "

$ws.Range("D9").Value = "2.428.73"
$ws.Range("E9").Value = "  -2.14%  "

$ws.Range("E10").Value = "  -5.23%  "

$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("E12").Value = "  -2.09%  "

$ws.Range("E13").Value = "  -2.99%  "

$ws.Range("E14").Value = "  -2.05%  "

$ws.Range("E15").Value = "  -6.07%  "

$ws.Range("D16").Value = "2.871.09"
$ws.Range("E16").Value = "  -1.67%  "

$ws.Range("D17").Value = "62.327.73"
$ws.Range("E17").Value = "  -1.50%  "

$ws.Range("D18").Value = "2.420.47"
$ws.Range("E18").Value = "  -2.40%  "

$ws.Range("E19").Value = "  -4.02%  "

$ws.Range("E20").Value = "  -3.12%  "

$ws.Range("D21").Value = "'324.84"
$ws.Range("E21").Value = "  -1.26%  "

$ws.Range("E22").Value = "  -2.63%  "

$ws.Range("D23").Value = "'2.01"
$ws.Range("E23").Value = "  +1.84%  "

$ws.Range("E24").Value = "  -4.07%  "

$ws.Range("D25").Value = "'65.17"
$ws.Range("E25").Value = "  -3.28%  "

$ws.Range("D26").Value = "'620.87"
$ws.Range("E26").Value = "  -1.52%  "

$ws.Range("E27").Value = "  +1.29%  "

$ws.Range("E28").Value = "  -9.46%  "

$ws.Range("D29").Value = "2.550.11"
$ws.Range("E29").Value = "  -1.76%  "

$ws.Range("D30").Value = "'0.973"
$ws.Range("E30").Value = "  -2.43%  "

$ws.Range("E31").Value = "  -4.96%  "

$ws.Range("E32").Value = "  -4.65%  "

$ws.Range("D33").Value = "'1.86"
$ws.Range("E33").Value = "  -3.59%  "

$ws.Range("E34").Value = "  -8.17%  "

$ws.Range("D35").Value = "'5.04"
$ws.Range("E35").Value = "  -3.29%  "

$ws.Range("E36").Value = "  +0.36%  "

$ws.Range("E37").Value = "  -6.60%  "

$ws.Range("E38").Value = "  -3.15%  "

$ws.Range("E39").Value = "  -2.30%  "

$ws.Range("D40").Value = "'147.03"
$ws.Range("E40").Value = "  +0.45%  "

$ws.Range("E41").Value = "  -5.76%  "

$ws.Range("E42").Value = "  -7.02%  "

$ws.Range("E43").Value = "  +1.39%  "

$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("E45").Value = "  -8.22%  "

$ws.Range("D46").Value = "'145.23"
$ws.Range("E46").Value = "  -3.58%  "

$ws.Range("E47").Value = "  -2.36%  "

$ws.Range("D48").Value = "'20.12"
$ws.Range("E48").Value = "  -4.67%  "

$ws.Range("E49").Value = "  -5.28%  "

$ws.Range("D50").Value = "'0.594"
$ws.Range("E50").Value = "  -2.86%  "

$ws.Range("E51").Value = "  -4.66%  "
